# Update outputs-r202: refresh the "quadratic-svm-score" sheet from the
# latest ful-path.csv export.
#
# 1) Row 1 header cells + the row-2 "Row" label are re-stamped with a
#    freshly (re-)applied text format, which — just like the previous
#    regeneration pass that produced style indices 3/4 — mints a new
#    paired set of cellXfs/border entries (text + date-time) rather than
#    reusing the old ones.
# 2) The numeric prediction score in B2 is refreshed with the new
#    computed value from the updated source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the text ("@") format to the header row + the row-2 label,
# touching each cell's border object (then clearing it back to "no
# line") so each one lands on a brand-new style entry instead of being
# silently deduped against the pre-existing style 3.
$headerCells = @("A1", "B1", "C1", "A2")
foreach ($addr in $headerCells) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Borders(5).LineStyle = -4142
}

# The regeneration pass also mints the matching date-time style variant
# (paired with the text variant above, mirroring styles 3/4) even though
# nothing on this sheet currently needs it; stamp it on a scratch cell
# and clear the cell so only the style entry remains behind.
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "m/d/yy h:mm"
$scratch.Borders(6).LineStyle = -4142
$scratch.Clear()

# Refresh the prediction score with the newly computed value.
$ws.Range("B2").Value = 16.462121324595831
